# Commit: "commit ui with new wiring connecting (#1090)"
#
# The workbook contains two lookup tables ("火灾报警" and "照明" sheets)
# that map block types to a wiring/circuit "type code" in columns C/D.
# This edit replaces the old short machine codes (E-FAS-WIRE, E-EFPS-WIRE,
# E-LITE-WIRE, ...) with human readable Chinese wiring descriptions.

$wb = $excel.ActiveWorkbook

# Map of old wire-type code -> new human readable Chinese description.
$wireMap = @{
    "E-FAS-WIRE"   = "火灾报警总线"
    "E-FAS-WIRE2"  = "总线+DC24V电源线"
    "E-FAS-WIRE3"  = "可燃气体探测系统布线"
    "E-FAS-WIRE4"  = "手动控制线"
    "E-FAS-WIRE5"  = "消防电话总线"
    "E-CTRL-WIRE"  = "消防联动控制线"
    "E-BRST-WIRE"  = "消防广播线"
    "E-FDS-WIRE"   = "防火门监控系统布线"
    "E-PMFE-WIRE"  = "消防电源监控系统布线"
    "E-EFPS-WIRE"  = "电气火灾监控系统布线"
    "E-LITE-WIRE"  = "平时照明布线"
    "E-LITE-WIRE2" = "消防应急照明布线"
}

# --- Sheet 1: 火灾报警 (fire alarm) ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$lastRow1 = 88
for ($r = 3; $r -le $lastRow1; $r++) {
    foreach ($col in 3, 4) {
        $cell = $ws1.Cells.Item($r, $col)
        $old = $cell.Value2
        if ($old -ne $null -and $wireMap.ContainsKey($old)) {
            $cell.Value2 = $wireMap[$old]
        }
    }
}

# Widen columns C/D now that they hold longer descriptive text.
$ws1.Columns.Item(3).ColumnWidth = 24
$ws1.Columns.Item(4).ColumnWidth = 21

# Move the saved selection/active cell back to the top of the sheet.
$ws1.Activate()
$ws1.Range("C2").Select()

# --- Sheet 2: 照明 (lighting) ----------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$lastRow2 = 86
for ($r = 3; $r -le $lastRow2; $r++) {
    $cell = $ws2.Cells.Item($r, 3)
    $old = $cell.Value2
    if ($old -ne $null -and $wireMap.ContainsKey($old)) {
        $cell.Value2 = $wireMap[$old]
    }
}
